$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host ("Sheet1 name: " + $ws.Name)
